# Generate Report for Handback
# Rewrites the Overview / zh-cn / de-de sheets so that:
#   - the previously-handed-off file (8a6eb94e-...) is now reported as
#     "Handed back: in sync with en-US" and its handoff/handback uuid+hash
#     pair is replaced by the new one (03a706f5-...)
#   - a second source file (ffffad353213-...) shows up with the same
#     handback status
#   - the .localization-config bookkeeping row moves down one row

$wb = $excel.ActiveWorkbook

$mdNew      = "03a706f5-8170-4e62-9d3e-fded9e2304b5.md"
$mdNew2     = "ffffad353213-e9eb-4881-8dfb-293c2a81b6bf.md"
$cfgName    = ".localization-config"

$statusHandback = "Handed back: in sync with en-US"
$statusNotLoc   = "Not to be localized"
$statusInclude  = "Include"
$statusIgnored  = "Ignored"

$xlfZh = "03a706f5-8170-4e62-9d3e-fded9e2304b5.baa2e4984f6cff7f5e57e681a81951a7580def8e.zh-cn.xlf"
$xlfDe = "03a706f5-8170-4e62-9d3e-fded9e2304b5.baa2e4984f6cff7f5e57e681a81951a7580def8e.de-de.xlf"

$dtZhHandoff  = "2016-02-29 04:22:04"
$dtZhHandback = "2016-02-29 04:22:47"
$dtDeHandoff  = "2016-02-29 04:22:14"
$dtDeHandback = "2016-02-29 04:23:06"
$dtZero       = "0001-01-01 00:00:00"

$baseRepo   = "https://github.com/OpenLocalizationTest/oltest/blob/504e37910e8be50e6b460e0124d22440086e24bc"
$baseHandZh = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/baa2e4984f6cff7f5e57e681a81951a7580def8e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$baseHandDe = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/baa2e4984f6cff7f5e57e681a81951a7580def8e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Hyperlinks.Delete()

$ov.Range("A2").Value = $mdNew
$ov.Range("B2").Value = $statusHandback
$ov.Range("C2").Value = $statusHandback

$ov.Range("A3").Value = $mdNew2
$ov.Range("B3").Value = $statusHandback
$ov.Range("C3").Value = $statusHandback

$ov.Range("A4").Value = $cfgName
$ov.Range("B4").Value = $statusNotLoc
$ov.Range("C4").Value = $statusNotLoc

$ov.Hyperlinks.Add($ov.Range("A2"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$ov.Hyperlinks.Add($ov.Range("A3"), "$baseRepo/e2e/$mdNew2", "", "", $mdNew2)
$ov.Hyperlinks.Add($ov.Range("A4"), "$baseRepo/$cfgName", "", "", $cfgName)

$ov.Range("A2").Style = "HyperLink"
$ov.Range("A3").Style = "HyperLink"
$ov.Range("A4").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Hyperlinks.Delete()

# row 2 - file that was just handed back
$zh.Range("A2").Value = $mdNew
$zh.Range("B2").Value = $statusHandback
$zh.Range("C2").Value = $xlfZh
$zh.Range("D2").Value = $dtZhHandoff
$zh.Range("E2").Value = $mdNew
$zh.Range("F2").Value = $xlfZh
$zh.Range("G2").Value = $dtZhHandback
$zh.Range("H2").Value = $statusInclude
$zh.Range("I2").Value = ""

# row 3 - second source file, same handback batch
$zh.Range("A3").Value = $mdNew2
$zh.Range("B3").Value = $statusHandback
$zh.Range("C3").Value = $xlfZh
$zh.Range("D3").Value = $dtZhHandoff
$zh.Range("E3").Value = $mdNew
$zh.Range("F3").Value = $xlfZh
$zh.Range("G3").Value = $dtZhHandback
$zh.Range("H3").Value = $statusInclude
$zh.Range("I3").Value = ""

# row 4 - bookkeeping row, shifted down from row 3
$zh.Range("A4").Value = $cfgName
$zh.Range("B4").Value = $statusNotLoc
$zh.Range("C4").Value = ""
$zh.Range("D4").Value = $dtZero
$zh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("E4").Value = ""
$zh.Range("F4").Value = ""
$zh.Range("G4").Value = $dtZero
$zh.Range("H4").Value = $statusIgnored

$zh.Hyperlinks.Add($zh.Range("A2"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$zh.Hyperlinks.Add($zh.Range("C2"), "$baseHandZh/ht/$xlfZh", "", "", $xlfZh)
$zh.Hyperlinks.Add($zh.Range("E2"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$zh.Hyperlinks.Add($zh.Range("F2"), "$baseHandZh/ht/$xlfZh", "", "", $xlfZh)
$zh.Hyperlinks.Add($zh.Range("A3"), "$baseRepo/e2e/$mdNew2", "", "", $mdNew2)
$zh.Hyperlinks.Add($zh.Range("C3"), "$baseHandZh/ht/$xlfZh", "", "", $xlfZh)
$zh.Hyperlinks.Add($zh.Range("E3"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$zh.Hyperlinks.Add($zh.Range("F3"), "$baseHandZh/ht/$xlfZh", "", "", $xlfZh)
$zh.Hyperlinks.Add($zh.Range("A4"), "$baseRepo/$cfgName", "", "", $cfgName)

$zh.Range("A2").Style = "HyperLink"
$zh.Range("C2").Style = "HyperLink"
$zh.Range("E2").Style = "HyperLink"
$zh.Range("F2").Style = "HyperLink"
$zh.Range("A3").Style = "HyperLink"
$zh.Range("C3").Style = "HyperLink"
$zh.Range("E3").Style = "HyperLink"
$zh.Range("F3").Style = "HyperLink"
$zh.Range("A4").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Cells.Hyperlinks.Delete()

# row 2 - file that was just handed back
$de.Range("A2").Value = $mdNew
$de.Range("B2").Value = $statusHandback
$de.Range("C2").Value = $xlfDe
$de.Range("D2").Value = $dtDeHandoff
$de.Range("E2").Value = $mdNew
$de.Range("F2").Value = $xlfDe
$de.Range("G2").Value = $dtDeHandback
$de.Range("H2").Value = $statusInclude
$de.Range("I2").Value = ""

# row 3 - second source file, same handback batch
$de.Range("A3").Value = $mdNew2
$de.Range("B3").Value = $statusHandback
$de.Range("C3").Value = $xlfDe
$de.Range("D3").Value = $dtDeHandoff
$de.Range("E3").Value = $mdNew
$de.Range("F3").Value = $xlfDe
$de.Range("G3").Value = $dtDeHandback
$de.Range("H3").Value = $statusInclude
$de.Range("I3").Value = ""

# row 4 - bookkeeping row, shifted down from row 3
$de.Range("A4").Value = $cfgName
$de.Range("B4").Value = $statusNotLoc
$de.Range("C4").Value = ""
$de.Range("D4").Value = $dtZero
$de.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("E4").Value = ""
$de.Range("F4").Value = ""
$de.Range("G4").Value = $dtZero
$de.Range("H4").Value = $statusIgnored

$de.Hyperlinks.Add($de.Range("A2"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$de.Hyperlinks.Add($de.Range("C2"), "$baseHandDe/ht/$xlfDe", "", "", $xlfDe)
$de.Hyperlinks.Add($de.Range("E2"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$de.Hyperlinks.Add($de.Range("F2"), "$baseHandDe/ht/$xlfDe", "", "", $xlfDe)
$de.Hyperlinks.Add($de.Range("A3"), "$baseRepo/e2e/$mdNew2", "", "", $mdNew2)
$de.Hyperlinks.Add($de.Range("C3"), "$baseHandDe/ht/$xlfDe", "", "", $xlfDe)
$de.Hyperlinks.Add($de.Range("E3"), "$baseRepo/e2e/$mdNew", "", "", $mdNew)
$de.Hyperlinks.Add($de.Range("F3"), "$baseHandDe/ht/$xlfDe", "", "", $xlfDe)
$de.Hyperlinks.Add($de.Range("A4"), "$baseRepo/$cfgName", "", "", $cfgName)

$de.Range("A2").Style = "HyperLink"
$de.Range("C2").Style = "HyperLink"
$de.Range("E2").Style = "HyperLink"
$de.Range("F2").Style = "HyperLink"
$de.Range("A3").Style = "HyperLink"
$de.Range("C3").Style = "HyperLink"
$de.Range("E3").Style = "HyperLink"
$de.Range("F3").Style = "HyperLink"
$de.Range("A4").Style = "HyperLink"

Write-Host "Handback report generated."
